# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# formatting of the existing header cells and filling in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (border/bold/alignment) from the existing header
# cell H1 onto the two new header cells so they match the rest of the
# header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (rows 2-19): I = 1, J = same value as column H ------------
for ($r = 2; $r -le 19; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# --- Row 20 is a special case: I20 = 7, J20 = 8 ---------------------------
$ws.Cells.Item(20, 9).Value = 7
$ws.Cells.Item(20, 10).Value = 8
